$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing data rows (2-9) before writing the new full dataset (rows 2-15)
$ws.Range("A2:I9").ClearContents()

$ws.Range("A2").Value = "87548754"
$ws.Range("B2").Value = "yanet altamirano quiroz"
$ws.Range("C2").Value = "2025-03-13"
$ws.Range("D2").Value = "07:55:53"
$ws.Range("E2").Value = "13:03:49"
$ws.Range("F2").Value = "14:11:23"
$ws.Range("G2").Value = "17:36:28"
$ws.Range("H2").Value = "0 minutos"
$ws.Range("I2").Value = ""

$ws.Range("A3").Value = "87548754"
$ws.Range("B3").Value = "yanet altamirano quiroz"
$ws.Range("C3").Value = "2025-03-15"
$ws.Range("D3").Value = "07:50:27"
$ws.Range("E3").Value = "13:05:18"
$ws.Range("F3").Value = "14:03:07"
$ws.Range("G3").Value = "18:44:33"
$ws.Range("H3").Value = "0 minutos"
$ws.Range("I3").Value = ""

$ws.Range("A4").Value = "12345667"
$ws.Range("B4").Value = "carla siares adrianzen"
$ws.Range("C4").Value = "2025-03-12"
$ws.Range("D4").Value = "07:55:22"
$ws.Range("E4").Value = "13:04:01"
$ws.Range("F4").ClearContents()
$ws.Range("G4").ClearContents()
$ws.Range("H4").Value = "0 minutos"
$ws.Range("I4").Value = "Tiene horas sin marcar"

$ws.Range("A5").Value = "87548754"
$ws.Range("B5").Value = "yanet altamirano quiroz"
$ws.Range("C5").Value = "2025-03-18"
$ws.Range("D5").Value = "07:53:38"
$ws.Range("E5").Value = "13:53:38"
$ws.Range("F5").ClearContents()
$ws.Range("G5").Value = "17:53:39"
$ws.Range("H5").Value = "0 minutos"
$ws.Range("I5").Value = "Tiene horas sin marcar"

$ws.Range("A6").Value = "87548754"
$ws.Range("B6").Value = "yanet altamirano quiroz"
$ws.Range("C6").Value = "2025-03-17"
$ws.Range("D6").Value = "07:52:27"
$ws.Range("E6").ClearContents()
$ws.Range("F6").Value = "14:14:52"
$ws.Range("G6").Value = "18:10:06"
$ws.Range("H6").Value = "0 minutos"
$ws.Range("I6").Value = "Tiene horas sin marcar"

$ws.Range("A7").Value = "87548754"
$ws.Range("B7").Value = "yanet altamirano quiroz"
$ws.Range("C7").Value = "2025-03-16"
$ws.Range("D7").Value = "07:51:52"
$ws.Range("E7").Value = "13:07:12"
$ws.Range("F7").Value = "14:06:03"
$ws.Range("G7").Value = "17:53:26"
$ws.Range("H7").Value = "0 minutos"
$ws.Range("I7").Value = ""

$ws.Range("A8").Value = "12345667"
$ws.Range("B8").Value = "carla siares adrianzen"
$ws.Range("C8").Value = "2025-03-16"
$ws.Range("D8").Value = "07:51:52"
$ws.Range("E8").Value = "13:07:12"
$ws.Range("F8").Value = "14:06:03"
$ws.Range("G8").Value = "17:53:26"
$ws.Range("H8").Value = "0 minutos"
$ws.Range("I8").Value = ""

$ws.Range("A9").Value = "87548754"
$ws.Range("B9").Value = "yanet altamirano quiroz"
$ws.Range("C9").Value = "2025-03-12"
$ws.Range("D9").Value = "07:55:22"
$ws.Range("E9").Value = "13:04:01"
$ws.Range("F9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").Value = "0 minutos"
$ws.Range("I9").Value = "Tiene horas sin marcar"

$ws.Range("A10").Value = "12345667"
$ws.Range("B10").Value = "carla siares adrianzen"
$ws.Range("C10").Value = "2025-03-17"
$ws.Range("D10").Value = "07:52:27"
$ws.Range("E10").ClearContents()
$ws.Range("F10").Value = "14:14:52"
$ws.Range("G10").Value = "18:10:06"
$ws.Range("H10").Value = "0 minutos"
$ws.Range("I10").Value = "Tiene horas sin marcar"

$ws.Range("A11").Value = "12345667"
$ws.Range("B11").Value = "carla siares adrianzen"
$ws.Range("C11").Value = "2025-03-18"
$ws.Range("D11").Value = "07:53:38"
$ws.Range("E11").ClearContents()
$ws.Range("F11").ClearContents()
$ws.Range("G11").ClearContents()
$ws.Range("H11").Value = "0 minutos"
$ws.Range("I11").Value = "Tiene horas sin marcar"

$ws.Range("A12").Value = "12345667"
$ws.Range("B12").Value = "carla siares adrianzen"
$ws.Range("C12").Value = "2025-03-15"
$ws.Range("D12").Value = "07:50:27"
$ws.Range("E12").Value = "13:05:18"
$ws.Range("F12").Value = "14:03:07"
$ws.Range("G12").Value = "18:44:33"
$ws.Range("H12").Value = "0 minutos"
$ws.Range("I12").Value = ""

$ws.Range("A13").Value = "12345667"
$ws.Range("B13").Value = "carla siares adrianzen"
$ws.Range("C13").Value = "2025-03-13"
$ws.Range("D13").Value = "07:55:53"
$ws.Range("E13").Value = "13:03:49"
$ws.Range("F13").Value = "14:11:23"
$ws.Range("G13").Value = "16:36:28"
$ws.Range("H13").Value = "24 minutos"
$ws.Range("I13").Value = ""

$ws.Range("A14").Value = "87548754"
$ws.Range("B14").Value = "yanet altamirano quiroz"
$ws.Range("C14").Value = "2025-03-14"
$ws.Range("D14").Value = "07:50:18"
$ws.Range("E14").Value = "13:31:46"
$ws.Range("F14").ClearContents()
$ws.Range("G14").ClearContents()
$ws.Range("H14").Value = "0 minutos"
$ws.Range("I14").Value = "Tiene horas sin marcar"

$ws.Range("A15").Value = "12345667"
$ws.Range("B15").Value = "carla siares adrianzen"
$ws.Range("C15").Value = "2025-03-14"
$ws.Range("D15").Value = "07:50:18"
$ws.Range("E15").Value = "13:31:46"
$ws.Range("F15").ClearContents()
$ws.Range("G15").ClearContents()
$ws.Range("H15").Value = "0 minutos"
$ws.Range("I15").Value = "Tiene horas sin marcar"
